$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first rule row (row 10: "Greater than USD50 CR") entirely,
# shifting all subsequent rows up by one.
$ws.Rows.Item(10).Delete()

# Match the resulting selection state (whole of the new row 10 selected).
[void]$ws.Rows.Item(10).Select()
